$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column R (year 2020) plus one value per existing data row (4-14),
# mirroring the formatting already applied to the corresponding cell in
# column Q (the previous last year of data).
$rValues = [ordered]@{
    4  = 2020
    5  = 5
    6  = 3.5
    7  = 1.8
    8  = 24.4
    9  = 7.2
    10 = 2.9
    11 = 7.4
    12 = 4
    13 = 3.2
    14 = 3.5
}

foreach ($row in $rValues.Keys) {
    $qCell = $ws.Range("Q$row")
    $rCell = $ws.Range("R$row")
    $qCell.Copy()
    $rCell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    $rCell.Value = $rValues[$row]
}
$excel.CutCopyMode = $false

# Update the sheet's active selection to match the new data range (column R,
# rows 4-14), as recorded for the updated workbook.
$ws.Range("R4:R14").Select()
